# F05 Froze Encoder 1
# Updates the ASR results table (columns A/B/C, rows 2-18) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = "<them>";    B = "<them>";    C = 9 }
    @{ Row = 3;  A = "<his>";     B = "<this>";    C = 8 }
    @{ Row = 4;  A = "<yankee>";  B = "<yankee>";  C = 4 }
    @{ Row = 5;  A = "<find>";    B = "<find>";    C = 15 }
    @{ Row = 6;  A = "<as>";      B = "<as>";      C = 9 }
    @{ Row = 7;  A = "<been>";    B = "<been>";    C = 5 }
    @{ Row = 8;  A = "<november>"; B = "<november>"; C = 10 }
    @{ Row = 9;  A = "<six>";     B = "<six>";     C = 5 }
    @{ Row = 10; A = "<control>"; B = "<control>"; C = 8 }
    @{ Row = 11; A = "<zero>";    B = "<zero>";    C = 8 }
    @{ Row = 12; A = "<other>";   B = "<other>";   C = 7 }
    @{ Row = 13; A = "<echo>";    B = "<echo>";    C = 10 }
    @{ Row = 14; A = "<five>";    B = "<five>";    C = 9 }
    @{ Row = 15; A = "<papa>";    B = "<papa>";    C = 7 }
    @{ Row = 16; A = "<look>";    B = "<look>";    C = 4 }
    @{ Row = 17; A = "<like>";    B = "<would>";   C = 4 }
    @{ Row = 18; A = "<yankee>";  B = "<yankee>";  C = 6 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
}
